$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas diarias")

# Row 6 ("Cambiar imagenes al visualizar los productos") moves from
# "Incompleto" to "Completo": copy the formatting used by the other
# "Completo" rows (e.g. row 2) onto row 6, then update the status text.
$ws.Range("A2:D2").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Completo"

# The row is shorter now that its content/format matches the "Completo"
# styling (no longer needs the taller wrap height).
$ws.Rows.Item(6).RowHeight = 30

# Update the view: scroll/selection moved to D5.
$ws.Activate()
$ws.Range("D5").Select()
